$wb = $excel.ActiveWorkbook

# ALC row 17 (Leve Item ID 38956)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 965.9091
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 965.9091
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2897.7273
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3233.7273

# ALC row 120 (Leve Item ID 26279)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 35700
$ws.Range("J120").Value = 35700
$ws.Range("L120").Value = 35700
$ws.Range("N120").Value = -45376

# ALC row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4603.852
$ws.Range("I132").Value = 4177.8945
$ws.Range("J132").Value = 5615.5
$ws.Range("K132").Value = 12533.6835
$ws.Range("L132").Value = 16846.5
$ws.Range("M132").Value = -10003.6835
$ws.Range("N132").Value = -21906.5

# ARM row 61 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2376.125
$ws.Range("I61").Value = 2865.8572
$ws.Range("J61").Value = 1995.2222
$ws.Range("K61").Value = 2865.8572
$ws.Range("L61").Value = 1995.2222
$ws.Range("M61").Value = -2653.8572
$ws.Range("N61").Value = -2419.2222

# ARM row 136 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2376.125
$ws.Range("I136").Value = 2865.8572
$ws.Range("J136").Value = 1995.2222
$ws.Range("K136").Value = 8597.571599999999
$ws.Range("L136").Value = 5985.6666
$ws.Range("M136").Value = -6047.571599999999
$ws.Range("N136").Value = -11085.6666

# BSM row 63 (Leve Item ID 10592)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51372

# BSM row 66 (Leve Item ID 10592)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156864

# BSM row 94 (Leve Item ID 19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2422.4167
$ws.Range("I94").Value = 2165.5715
$ws.Range("J94").Value = 2782
$ws.Range("K94").Value = 2165.5715
$ws.Range("L94").Value = 2782
$ws.Range("M94").Value = -1714.5715
$ws.Range("N94").Value = -3684

# BSM row 132 (Leve Item ID 41855)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 29900
$ws.Range("J132").Value = 29900
$ws.Range("L132").Value = 29900
$ws.Range("N132").Value = -40020

# CRP row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3625065.5
$ws.Range("I31").Value = 1198.8032
$ws.Range("J31").Value = 10755900
$ws.Range("K31").Value = 1198.8032
$ws.Range("L31").Value = 10755900
$ws.Range("M31").Value = -903.8032000000001
$ws.Range("N31").Value = -10756490

# CRP row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3625065.5
$ws.Range("I34").Value = 1198.8032
$ws.Range("J34").Value = 10755900
$ws.Range("K34").Value = 1198.8032
$ws.Range("L34").Value = 10755900
$ws.Range("M34").Value = -996.8032000000001
$ws.Range("N34").Value = -10756304

# CRP row 43 (Leve Item ID 18504)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 18197.545
$ws.Range("J43").Value = 18197.545
$ws.Range("L43").Value = 18197.545
$ws.Range("N43").Value = -18565.545

# CRP row 70 (Leve Item ID 12011)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 12666.667
$ws.Range("J70").Value = 12666.667
$ws.Range("L70").Value = 12666.667
$ws.Range("N70").Value = -13296.667

# CRP row 73 (Leve Item ID 12011)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 12666.667
$ws.Range("J73").Value = 12666.667
$ws.Range("L73").Value = 12666.667
$ws.Range("N73").Value = -14850.667

# CRP row 101 (Leve Item ID 18504)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H101").Value = 18197.545
$ws.Range("J101").Value = 18197.545
$ws.Range("L101").Value = 18197.545
$ws.Range("N101").Value = -24687.545

# CUL row 11 (Leve Item ID 4745)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 528.6
$ws.Range("I11").Value = 28.666666
$ws.Range("J11").Value = 742.8570999999999
$ws.Range("K11").Value = 85.99999800000001
$ws.Range("L11").Value = 2228.5713
$ws.Range("M11").Value = 54.00000199999999
$ws.Range("N11").Value = -2508.5713

# CUL row 26 (Leve Item ID 4746)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 351
$ws.Range("I26").Value = 107.875
$ws.Range("J26").Value = 740
$ws.Range("K26").Value = 323.625
$ws.Range("L26").Value = 2220
$ws.Range("M26").Value = -35.625
$ws.Range("N26").Value = -2796

# CUL row 52 (Leve Item ID 31902)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1461.625
$ws.Range("J52").Value = 1461.625
$ws.Range("L52").Value = 4384.875
$ws.Range("N52").Value = -4916.875

# CUL row 93 (Leve Item ID 19808)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 6398.5
$ws.Range("J93").Value = 6398.5
$ws.Range("L93").Value = 19195.5
$ws.Range("N93").Value = -22939.5

# CUL row 109 (Leve Item ID 27854)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2128.125
$ws.Range("I109").Value = 1326.6666
$ws.Range("J109").Value = 3158.5715
$ws.Range("K109").Value = 3979.9998
$ws.Range("L109").Value = 9475.7145
$ws.Range("M109").Value = -2939.9998
$ws.Range("N109").Value = -11555.7145

# CUL row 115 (Leve Item ID 27861)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 4105.4
$ws.Range("I115").Value = 3464
$ws.Range("J115").Value = 4533
$ws.Range("K115").Value = 10392
$ws.Range("L115").Value = 13599
$ws.Range("M115").Value = -9217
$ws.Range("N115").Value = -15949

# CUL row 124 (Leve Item ID 36040)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 8171
$ws.Range("I124").Value = 4030
$ws.Range("J124").Value = 8999.200000000001
$ws.Range("K124").Value = 12090
$ws.Range("L124").Value = 26997.6
$ws.Range("M124").Value = -7180
$ws.Range("N124").Value = -36817.60000000001

# CUL row 125 (Leve Item ID 36043)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 7875
$ws.Range("J125").Value = 7875
$ws.Range("L125").Value = 23625
$ws.Range("N125").Value = -33465

# CUL row 140 (Leve Item ID 44097)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1279.6774
$ws.Range("I140").Value = 920.55554
$ws.Range("J140").Value = 1776.9231
$ws.Range("K140").Value = 2761.66662
$ws.Range("L140").Value = 5330.7693
$ws.Range("M140").Value = 2418.33338
$ws.Range("N140").Value = -15690.7693

# GSM row 92 (Leve Item ID 18094)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 5833.4443
$ws.Range("J92").Value = 5833.4443
$ws.Range("L92").Value = 5833.4443
$ws.Range("N92").Value = -9577.444299999999

# LTW row 16 (Leve Item ID 5289)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3319.65
$ws.Range("I16").Value = 1887.3572
$ws.Range("J16").Value = 6661.6665
$ws.Range("K16").Value = 1887.3572
$ws.Range("L16").Value = 6661.6665
$ws.Range("M16").Value = -1717.3572
$ws.Range("N16").Value = -7001.6665

# LTW row 74 (Leve Item ID 11990)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# LTW row 77 (Leve Item ID 11990)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# WVR row 33 (Leve Item ID 2734)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

# WVR row 36 (Leve Item ID 2734)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

# WVR row 126 (Leve Item ID 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3727.4092
$ws.Range("I126").Value = 3764.0557
$ws.Range("K126").Value = 11292.1671
$ws.Range("M126").Value = -8822.167099999999
